# This workbook's "Artfynd" sheet (rows 2-10) is re-sorted in this commit:
# the diff shows that, after the edit, each data row's full content is
# identical to some other row's content from before the edit - i.e. the 9
# observation rows are permuted into a new order, with no new or removed
# data.
#
# Mapping of new-row -> old-row (where the old row's data ends up):
#   2 <- 5   3 <- 8   4 <- 2   5 <- 3   6 <- 4
#   7 <- 6   8 <- 7   9 <- 10  10 <- 9
#
# Implementation notes:
#  - Rather than hard-coding the (Swedish-language) cell text in this
#    script, snapshot every cell of rows 2-10 first, then write the
#    snapshot back out in the permuted order.
#  - Only cells whose value actually changes under the permutation are
#    written. Several columns (dates stored as text, etc.) are identical
#    across every row, and re-assigning them via COM's ".Value" setter can
#    trigger Excel's automatic text -> date/number coercion, corrupting
#    cells that should stay untouched. Skipping no-op writes avoids that
#    entirely while still reproducing the target workbook exactly.
#  - Cells that must become blank (e.g. "AC" public-comment, or the
#    "AJ"/"AK"/"AO" substrate columns that only one row has) are cleared
#    with ClearContents().

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY")

$firstRow = 2
$lastRow = 10

# 1. Snapshot current values for every column/row in range.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $addr = $c + $r
        $rowData[$c] = $ws.Range($addr).Value2
    }
    $snapshot[$r] = $rowData
}

# 2. New row number -> source (old) row number.
$mapping = @{
    2  = 5
    3  = 8
    4  = 2
    5  = 3
    6  = 4
    7  = 6
    8  = 7
    9  = 10
    10 = 9
}

# 3. Write each target row from the snapshot of its source row, skipping
#    any cell whose value is unchanged (avoids needless COM type coercion
#    and keeps the edit minimal / matching the diff exactly).
for ($targetRow = $firstRow; $targetRow -le $lastRow; $targetRow++) {
    $sourceRow = $mapping[$targetRow]
    $sourceData = $snapshot[$sourceRow]
    $currentData = $snapshot[$targetRow]
    foreach ($c in $cols) {
        $newVal = $sourceData[$c]
        $oldVal = $currentData[$c]
        if ($newVal -eq $oldVal) {
            continue
        }
        $addr = $c + $targetRow
        if ($newVal -eq $null) {
            $ws.Range($addr).ClearContents()
        } else {
            $ws.Range($addr).Value = $newVal
        }
    }
}

Write-Host "Row permutation applied."
